# Add scenario files for runs March 8th 2023
# - Clear the H2 coverage value on the "Platform Coverage" sheet
# - Add a new "Vector Control" platform row (row 12) with 0.25 coverage
#   values on every other year column from 2026 through 2040

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Remove the stray coverage value previously stored in H2
$ws.Range("H2").ClearContents()

# New row describing the "Vector Control" platform
$ws.Range("B12").Value = "Vector Control"

$ws.Range("X12").Value = 0.25
$ws.Range("Z12").Value = 0.25
$ws.Range("AB12").Value = 0.25
$ws.Range("AD12").Value = 0.25
$ws.Range("AF12").Value = 0.25
$ws.Range("AH12").Value = 0.25
$ws.Range("AJ12").Value = 0.25
$ws.Range("AL12").Value = 0.25
$ws.Range("AN12").Value = 0.25
$ws.Range("AP12").Value = 0.25
$ws.Range("AR12").Value = 0.25
$ws.Range("AT12").Value = 0.25
$ws.Range("AV12").Value = 0.25
$ws.Range("AX12").Value = 0.25
$ws.Range("AZ12").Value = 0.25

Write-Output "Applied Vector Control scenario row and cleared H2"
